$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36..129 down to 37..130
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record
$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 44497
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = 100112039
$ws.Cells.Item(36, 7).Value = "Ciboulette"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 80
$ws.Cells.Item(36, 11).Value = 2500
$ws.Cells.Item(36, 12).Value = 3000
$ws.Cells.Item(36, 13).Value = 2750
$ws.Cells.Item(36, 14).Value = "$/docena de atados"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 917
$ws.Cells.Item(36, 17).Value = 3
$ws.Cells.Item(36, 18).Value = "Hortaliza"
